$wb = $excel.ActiveWorkbook

# Hunk 0: ALC!row123
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Hunk 1: ALC!row129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 268075.38
$ws.Range("I129").Value = 4003492.2
$ws.Range("J129").Value = 1259.881
$ws.Range("K129").Value = 12010476.6
$ws.Range("L129").Value = 3779.643
$ws.Range("M129").Value = -12005476.6
$ws.Range("N129").Value = -13779.643

# Hunk 2: ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16427.4
$ws.Range("I32").Value = 15600.786
$ws.Range("J32").Value = 28000
$ws.Range("K32").Value = 15600.786
$ws.Range("L32").Value = 28000
$ws.Range("M32").Value = -15313.786
$ws.Range("N32").Value = -28574

# Hunk 3: ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 15153894
$ws.Range("I132").Value = 21740974
$ws.Range("J132").Value = 3609.4
$ws.Range("K132").Value = 65222922
$ws.Range("L132").Value = 10828.2
$ws.Range("M132").Value = -65220392
$ws.Range("N132").Value = -15888.2

# Hunk 4: BSM!row105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1972.94
$ws.Range("I105").Value = 1788.3928
$ws.Range("J105").Value = 2207.818
$ws.Range("K105").Value = 1788.3928
$ws.Range("L105").Value = 2207.818
$ws.Range("M105").Value = -41.39280000000008
$ws.Range("N105").Value = -5701.818

# Hunk 5: BSM!row132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 58747.5
$ws.Range("J132").Value = 58747.5
$ws.Range("L132").Value = 58747.5
$ws.Range("N132").Value = -68867.5

# Hunk 6: BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1907.228
$ws.Range("I134").Value = 1689
$ws.Range("J134").Value = 2090.258
$ws.Range("K134").Value = 5067
$ws.Range("L134").Value = 6270.773999999999
$ws.Range("M134").Value = -2532
$ws.Range("N134").Value = -11340.774

# Hunk 7: CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6191.9062
$ws.Range("I31").Value = 4001.8572
$ws.Range("J31").Value = 6805.12
$ws.Range("K31").Value = 4001.8572
$ws.Range("L31").Value = 6805.12
$ws.Range("M31").Value = -3706.8572
$ws.Range("N31").Value = -7395.12

# Hunk 8: CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6191.9062
$ws.Range("I34").Value = 4001.8572
$ws.Range("J34").Value = 6805.12
$ws.Range("K34").Value = 4001.8572
$ws.Range("L34").Value = 6805.12
$ws.Range("M34").Value = -3799.8572
$ws.Range("N34").Value = -7209.12

# Hunk 9: CRP!row59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 37804.875
$ws.Range("J59").Value = 38491.285
$ws.Range("L59").Value = 38491.285
$ws.Range("N59").Value = -40781.285

# Hunk 10: CRP!row107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 722.0833
$ws.Range("I107").Value = 572
$ws.Range("J107").Value = 932.2
$ws.Range("K107").Value = 572
$ws.Range("L107").Value = 932.2
$ws.Range("M107").Value = 1348
$ws.Range("N107").Value = -4772.2

# Hunk 11: CRP!row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 166870560
$ws.Range("I122").Value = 250300600
$ws.Range("J122").Value = 10457
$ws.Range("K122").Value = 750901800
$ws.Range("L122").Value = 31371
$ws.Range("M122").Value = -750899350
$ws.Range("N122").Value = -36271

# Hunk 12: CUL!row134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 50202588
$ws.Range("I134").Value = 55779984
$ws.Range("J134").Value = 6033
$ws.Range("K134").Value = 167339952
$ws.Range("L134").Value = 18099
$ws.Range("M134").Value = -167334882
$ws.Range("N134").Value = -28239

# Hunk 13: GSM!row107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 89148.35000000001
$ws.Range("I107").Value = 250356.5
$ws.Range("J107").Value = 3170.6667
$ws.Range("K107").Value = 250356.5
$ws.Range("L107").Value = 3170.6667
$ws.Range("M107").Value = -248436.5
$ws.Range("N107").Value = -7010.6667

# Hunk 14: GSM!row126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3487.2856
$ws.Range("I126").Value = 4002.2
$ws.Range("J126").Value = 2200
$ws.Range("K126").Value = 12006.6
$ws.Range("L126").Value = 6600
$ws.Range("M126").Value = -9536.599999999999
$ws.Range("N126").Value = -11540

# Hunk 15: GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4330.1816
$ws.Range("I132").Value = 4396.0415
$ws.Range("J132").Value = 4154.5557
$ws.Range("K132").Value = 13188.1245
$ws.Range("L132").Value = 12463.6671
$ws.Range("M132").Value = -10658.1245
$ws.Range("N132").Value = -17523.6671

# Hunk 16: LTW!row68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2333.3333
$ws.Range("I68").Value = 2300
$ws.Range("J68").Value = 2500
$ws.Range("K68").Value = 2300
$ws.Range("L68").Value = 2500
$ws.Range("M68").Value = -1551
$ws.Range("N68").Value = -3998

# Hunk 17: LTW!row71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2333.3333
$ws.Range("I71").Value = 2300
$ws.Range("J71").Value = 2500
$ws.Range("K71").Value = 11500
$ws.Range("L71").Value = 12500
$ws.Range("M71").Value = -7756
$ws.Range("N71").Value = -19988

# Hunk 18: LTW!row122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 35531.035
$ws.Range("I122").Value = 43413.793
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 130241.379
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -127791.379
$ws.Range("N122").Value = -16900

# Hunk 19: LTW!row123
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 42359.285
$ws.Range("J123").Value = 42359.285
$ws.Range("L123").Value = 42359.285
$ws.Range("N123").Value = -52159.285

# Hunk 20: LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2977.7222
$ws.Range("I132").Value = 2654.6135
$ws.Range("J132").Value = 4399.4
$ws.Range("K132").Value = 7963.8405
$ws.Range("L132").Value = 13198.2
$ws.Range("M132").Value = -5433.8405
$ws.Range("N132").Value = -18258.2

# Hunk 21: WVR!row69
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

# Hunk 22: WVR!row72
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

# Hunk 23: WVR!row81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1277.1428
$ws.Range("I81").Value = 1277.1428
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2554.2856
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1493.2856
$ws.Range("N81").ClearContents()

# Hunk 24: WVR!row84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1277.1428
$ws.Range("I84").Value = 1277.1428
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 12771.428
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -7467.428

# Hunk 25: WVR!row98
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

# Hunk 26: WVR!row104
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# Hunk 27: WVR!row107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 16666913
$ws.Range("I107").Value = 269.75
$ws.Range("J107").Value = 50000200
$ws.Range("K107").Value = 809.25
$ws.Range("L107").Value = 150000600
$ws.Range("M107").Value = 1110.75
$ws.Range("N107").Value = -150004440

# Hunk 28: WVR!row108
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# Hunk 29: WVR!row113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 614.0833
$ws.Range("I113").Value = 579
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1737
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 433
$ws.Range("N113").Value = -7340

# Hunk 30: WVR!row114
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# Hunk 31: WVR!row122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 52859044
$ws.Range("J122").Value = 3249.75
$ws.Range("L122").Value = 9749.25
$ws.Range("N122").Value = -14649.25
